$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AA2").Value = 'maa://21246 (91.2), maa://36684 (98.57), ***maa://22731 (6.67)'
$ws.Range("AE2").Value = 'maa://25251 (92.41), ***maa://21730 (17.19), ***maa://39501 (25.0), *maa://36675 (60.0)'
$ws.Range("K3").Value = '*maa://22880 (69.74), maa://20276 (82.73), *maa://22749 (62.5)'
$ws.Range("W4").Value = '**maa://32495 (47.93), ***maa://31785 (18.02), ***maa://36683 (26.67)'
$ws.Range("AA4").Value = '*maa://32658 (73.33)'
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = '1'
$ws.Range("C6").Value = '**maa://42407 (50.0)'
$ws.Range("AA6").Value = 'maa://22739 (91.3)'
$ws.Range("AE6").Value = '*maa://33152 (59.38), ***maa://22770 (28.57)'
$ws.Range("K7").Value = 'maa://28624 (91.3), maa://24957 (97.37)'
$ws.Range("S7").Value = 'maa://21291 (89.19)'
$ws.Range("C8").Value = '*maa://21476 (69.77), **maa://39431 (40.0), *maa://37551 (57.14)'
$ws.Range("O8").Value = 'maa://32931 (88.46), *maa://21916 (60.34), maa://23252 (92.31), **maa://22759 (45.45), maa://37496 (100.0)'
$ws.Range("W8").Value = 'maa://21411 (96.0)'
$ws.Range("AE8").Value = '*maa://24479 (76.39), *maa://21990 (53.85)'
$ws.Range("K9").Value = 'maa://22762 (91.57), *maa://39552 (66.67)'
$ws.Range("W9").Value = 'maa://26223 (96.91)'
$ws.Range("AA9").Value = 'maa://28711 (87.95), ***maa://22740 (5.88), **maa://27377 (46.15), ***maa://25174 (20.0), **maa://39938 (50.0), maa://40166 (100.0)'
$ws.Range("O10").Value = 'maa://28977 (94.67), *maa://23264 (62.96), maa://36669 (91.3)'
$ws.Range("S10").Value = 'maa://27395 (96.67), maa://22755 (87.62), **maa://22756 (40.91), ***maa://21737 (10.61)'
$ws.Range("W10").Value = 'maa://22301 (97.4), maa://22726 (100.0)'
$ws.Range("S11").Value = 'maa://22747 (94.33), maa://22501 (98.15)'
$ws.Range("AA12").Value = 'maa://23669 (95.86), maa://36677 (94.87), maa://39872 (84.62)'
$ws.Range("AE12").Value = '*maa://28932 (78.63), *maa://20106 (63.64), *maa://22769 (62.96)'
$ws.Range("C13").Value = 'maa://24999 (91.42), maa://36673 (91.8), maa://25001 (85.51)'
$ws.Range("G13").Value = '*maa://21248 (75.61), **maa://22728 (47.62)'
$ws.Range("O13").Value = 'maa://22676 (91.84), *maa://22583 (75.0), *maa://22500 (55.81)'
$ws.Range("C14").Value = 'maa://30764 (86.05)'
$ws.Range("AE15").Value = 'maa://21364 (80.61), *maa://22766 (73.0), *maa://36666 (77.42)'
$ws.Range("C16").Value = 'maa://21441 (96.17), maa://36679 (93.75), maa://37650 (95.45)'
$ws.Range("S16").Value = 'maa://22729 (95.17), *maa://28648 (69.09), *maa://36674 (79.31)'
$ws.Range("G17").Value = 'maa://22430 (88.57), maa://39599 (84.21)'
$ws.Range("C18").Value = 'maa://24570 (96.51)'
$ws.Range("K18").Value = 'maa://22466 (88.46), *maa://22732 (51.85)'
$ws.Range("AA19").Value = '*maa://30709 (60.59), *maa://36668 (52.17)'
$ws.Range("G20").Value = 'maa://22864 (88.46)'
$ws.Range("K20").Value = 'maa://41331 (90.91)'
$ws.Range("W21").Value = 'maa://20110 (86.57), maa://34946 (90.62)'
$ws.Range("AA21").Value = '*maa://21443 (78.96), **maa://23820 (30.91)'
$ws.Range("AE21").Value = 'maa://22524 (94.29), *maa://22432 (74.55)'
$ws.Range("G22").Value = 'maa://25236 (95.83), **maa://21678 (48.94), **maa://22735 (50.0)'
$ws.Range("K22").Value = 'maa://27127 (86.52), *maa://22751 (77.05)'
$ws.Range("W22").Value = 'maa://21282 (98.82), *maa://37649 (66.67)'
$ws.Range("K23").Value = 'maa://39756 (92.21), maa://39875 (95.65)'
$ws.Range("O23").Value = 'maa://30587 (91.62), *maa://29748 (75.2), ***maa://29785 (15.15), *maa://37566 (78.95)'
$ws.Range("C24").Value = 'maa://24368 (80.73)'
$ws.Range("W24").Value = 'maa://23504 (92.88), maa://29988 (86.07), **maa://22892 (40.14), *maa://25141 (76.86), maa://36663 (80.36), ***maa://22815 (23.08)'
$ws.Range("AE24").Value = 'maa://22523 (84.86), *maa://36672 (76.74), maa://29910 (94.0), **maa://21440 (34.55)'
$ws.Range("G25").Value = '*maa://29063 (76.15), *maa://25311 (74.19), ***maa://22725 (4.84)'
$ws.Range("AA26").Value = '*maa://42235 (66.67)'
$ws.Range("AE26").Value = 'maa://30511 (84.38), *maa://29760 (61.54)'
$ws.Range("G27").Value = '**maa://21283 (48.65), maa://34494 (100.0), **maa://36665 (44.44), maa://39601 (88.89)'
$ws.Range("C28").Value = 'maa://24465 (90.33), maa://25725 (82.28)'
$ws.Range("S28").Value = 'maa://23263 (94.32), *maa://29765 (61.19)'
$ws.Range("W28").Value = 'maa://39929 (85.79), ***maa://39723 (14.71), maa://41749 (86.67)'
$ws.Range("AE29").Value = '*maa://24080 (68.87), ***maa://34960 (9.09)'
$ws.Range("K30").Value = 'maa://30442 (94.34)'
$ws.Range("O30").Value = 'maa://21442 (99.49)'
$ws.Range("K31").Value = 'maa://35926 (93.56), maa://36258 (80.52)'
$ws.Range("S32").Value = 'maa://41108 (90.32), maa://41238 (94.44)'
$ws.Range("AE32").ClearContents()
$ws.Range("O33").Value = '*maa://21956 (79.07), maa://22730 (82.14)'
$ws.Range("K35").Value = 'maa://41296 (97.73)'
$ws.Range("O38").Value = '*maa://24383 (66.67)'
$ws.Range("G39").Value = 'maa://25199 (86.11), maa://36670 (89.39), maa://30434 (87.27), ***maa://25036 (16.0)'
$ws.Range("O40").Value = 'maa://23278 (95.88), maa://21386 (95.63), maa://36664 (90.24)'
$ws.Range("O41").Value = '**maa://35616 (37.04)'
$ws.Range("G44").Value = 'maa://29768 (97.52), maa://27728 (96.0)'
$ws.Range("G45").Value = 'maa://21229 (85.47), maa://30807 (95.08), *maa://22767 (52.94), ***maa://20796 (13.79)'
$ws.Range("G46").Value = 'maa://35931 (92.48)'
$ws.Range("G47").Value = 'maa://27410 (95.81), maa://29661 (97.64), maa://28038 (84.62)'
$ws.Range("G53").Value = 'maa://32534 (93.21), **maa://32434 (36.36)'
$ws.Range("G55").Value = 'maa://32532 (92.31)'
$ws.Range("G59").Value = 'maa://27746 (84.0), maa://31270 (96.97)'
